# Auto-generated edit script: update Leve profit-calculation values
# across multiple sheets, per the authoritative diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 880
$ws.Range("I18").Value = 868
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 868
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = -584
$ws.Range("N18").Value = -1568
$ws.Range("H40").Value = 1715.3846
$ws.Range("I40").Value = 1542.8572
$ws.Range("J40").Value = 1916.6666
$ws.Range("K40").Value = 1542.8572
$ws.Range("L40").Value = 1916.6666
$ws.Range("M40").Value = -1367.8572
$ws.Range("N40").Value = -2266.6666
$ws.Range("H98").Value = 1464.8846
$ws.Range("I98").Value = 1153.75
$ws.Range("J98").Value = 2502
$ws.Range("K98").Value = 1153.75
$ws.Range("L98").Value = 2502
$ws.Range("M98").Value = 344.25
$ws.Range("N98").Value = -5498
$ws.Range("H122").Value = 1464.8846
$ws.Range("I122").Value = 1153.75
$ws.Range("J122").Value = 2502
$ws.Range("K122").Value = 3461.25
$ws.Range("L122").Value = 7506
$ws.Range("M122").Value = -1011.25
$ws.Range("N122").Value = -12406
$ws.Range("H137").Value = 329620.56
$ws.Range("I137").Value = 616770.4399999999
$ws.Range("J137").Value = 3313.9092
$ws.Range("K137").Value = 1850311.32
$ws.Range("L137").Value = 9941.7276
$ws.Range("M137").Value = -1847761.32
$ws.Range("N137").Value = -15041.7276
$ws.Range("H141").Value = 2393.5833
$ws.Range("I141").Value = 1901.5834
$ws.Range("J141").Value = 2639.5833
$ws.Range("K141").Value = 5704.7502
$ws.Range("L141").Value = 7918.749899999999
$ws.Range("M141").Value = -524.7502000000004
$ws.Range("N141").Value = -18278.7499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("H32").Value = 20384.455
$ws.Range("I32").Value = 6080.0347
$ws.Range("K32").Value = 6080.0347
$ws.Range("M32").Value = -5793.0347
$ws.Range("H61").Value = 2543.7021
$ws.Range("I61").Value = 2409.1628
$ws.Range("J61").Value = 3990
$ws.Range("K61").Value = 2409.1628
$ws.Range("L61").Value = 3990
$ws.Range("M61").Value = -2197.1628
$ws.Range("N61").Value = -4414
$ws.Range("H74").Value = 2050.2122
$ws.Range("I74").Value = 1175.5667
$ws.Range("K74").Value = 1175.5667
$ws.Range("M74").Value = -301.5667000000001
$ws.Range("H77").Value = 2050.2122
$ws.Range("I77").Value = 1175.5667
$ws.Range("K77").Value = 5877.833500000001
$ws.Range("M77").Value = -1509.833500000001
$ws.Range("H102").Value = 90910664
$ws.Range("I102").Value = 1630
$ws.Range("J102").Value = 250001470
$ws.Range("K102").Value = 1630
$ws.Range("L102").Value = 250001470
$ws.Range("M102").Value = -8
$ws.Range("N102").Value = -250004714
$ws.Range("H114").Value = 54000
$ws.Range("J114").Value = 54000
$ws.Range("L114").Value = 54000
$ws.Range("N114").Value = -62678
$ws.Range("H132").Value = 2232.4
$ws.Range("I132").Value = 2039.8
$ws.Range("J132").Value = 2617.6
$ws.Range("K132").Value = 6119.4
$ws.Range("L132").Value = 7852.799999999999
$ws.Range("M132").Value = -3589.4
$ws.Range("N132").Value = -12912.8
$ws.Range("H136").Value = 2543.7021
$ws.Range("I136").Value = 2409.1628
$ws.Range("J136").Value = 3990
$ws.Range("K136").Value = 7227.4884
$ws.Range("L136").Value = 11970
$ws.Range("M136").Value = -4677.4884
$ws.Range("N136").Value = -17070
$ws.Range("M26").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3007.651
$ws.Range("I31").Value = 1373.2972
$ws.Range("J31").Value = 5333.4614
$ws.Range("K31").Value = 1373.2972
$ws.Range("L31").Value = 5333.4614
$ws.Range("M31").Value = -1078.2972
$ws.Range("N31").Value = -5923.4614
$ws.Range("H34").Value = 3007.651
$ws.Range("I34").Value = 1373.2972
$ws.Range("J34").Value = 5333.4614
$ws.Range("K34").Value = 1373.2972
$ws.Range("L34").Value = 5333.4614
$ws.Range("M34").Value = -1171.2972
$ws.Range("N34").Value = -5737.4614
$ws.Range("H99").Value = 1917.6666
$ws.Range("I99").Value = 1512
$ws.Range("J99").Value = 1998.8
$ws.Range("K99").Value = 1512
$ws.Range("L99").Value = 1998.8
$ws.Range("M99").Value = -14
$ws.Range("N99").Value = -4994.8
$ws.Range("H126").Value = 1917.6666
$ws.Range("I126").Value = 1512
$ws.Range("J126").Value = 1998.8
$ws.Range("K126").Value = 4536
$ws.Range("L126").Value = 5996.4
$ws.Range("M126").Value = -2066
$ws.Range("N126").Value = -10936.4
$ws.Range("H132").Value = 2348.1516
$ws.Range("I132").Value = 1774
$ws.Range("J132").Value = 3879.2222
$ws.Range("K132").Value = 5322
$ws.Range("L132").Value = 11637.6666
$ws.Range("M132").Value = -2792
$ws.Range("N132").Value = -16697.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 3812.5
$ws.Range("I117").Value = 500
$ws.Range("J117").Value = 4285.7144
$ws.Range("K117").Value = 1500
$ws.Range("L117").Value = 12857.1432
$ws.Range("M117").Value = 1942
$ws.Range("N117").Value = -19741.1432
$ws.Range("H129").Value = 1468.8462
$ws.Range("I129").Value = 577.6923
$ws.Range("J129").Value = 2360
$ws.Range("K129").Value = 1733.0769
$ws.Range("L129").Value = 7080
$ws.Range("M129").Value = 3266.9231
$ws.Range("N129").Value = -17080
$ws.Range("H131").Value = 877.0612
$ws.Range("I131").Value = 462.05884
$ws.Range("J131").Value = 964.1605
$ws.Range("K131").Value = 1386.17652
$ws.Range("L131").Value = 2892.4815
$ws.Range("M131").Value = 3653.82348
$ws.Range("N131").Value = -12972.4815

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 28441.582
$ws.Range("I70").Value = 35696.97
$ws.Range("J70").Value = 4498.8
$ws.Range("K70").Value = 35696.97
$ws.Range("L70").Value = 4498.8
$ws.Range("M70").Value = -35426.97
$ws.Range("N70").Value = -5038.8
$ws.Range("H73").Value = 28441.582
$ws.Range("I73").Value = 35696.97
$ws.Range("J73").Value = 4498.8
$ws.Range("K73").Value = 35696.97
$ws.Range("L73").Value = 4498.8
$ws.Range("M73").Value = -34760.97
$ws.Range("N73").Value = -6370.8
$ws.Range("H103").Value = 37222
$ws.Range("J103").Value = 37222
$ws.Range("L103").Value = 37222
$ws.Range("N103").Value = -39566
$ws.Range("H122").Value = 2202.4075
$ws.Range("I122").Value = 2166.2273
$ws.Range("J122").Value = 2361.6
$ws.Range("K122").Value = 6498.6819
$ws.Range("L122").Value = 7084.799999999999
$ws.Range("M122").Value = -4048.6819
$ws.Range("N122").Value = -11984.8
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H132").Value = 2134.878
$ws.Range("I132").Value = 1952.2258
$ws.Range("J132").Value = 2701.1
$ws.Range("K132").Value = 5856.6774
$ws.Range("L132").Value = 8103.299999999999
$ws.Range("M132").Value = -3326.6774
$ws.Range("N132").Value = -13163.3
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4212.2666
$ws.Range("I132").Value = 5071.3335
$ws.Range("J132").Value = 2923.6667
$ws.Range("K132").Value = 15214.0005
$ws.Range("L132").Value = 8771.000100000001
$ws.Range("M132").Value = -12684.0005
$ws.Range("N132").Value = -13831.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8739.968999999999
$ws.Range("I122").Value = 11356.381
$ws.Range("J122").Value = 3745
$ws.Range("K122").Value = 34069.143
$ws.Range("L122").Value = 11235
$ws.Range("M122").Value = -31619.143
$ws.Range("N122").Value = -16135
$ws.Range("H132").Value = 2133.9524
$ws.Range("I132").Value = 2273.2666
$ws.Range("J132").Value = 1785.6666
$ws.Range("K132").Value = 6819.7998
$ws.Range("L132").Value = 5356.9998
$ws.Range("M132").Value = -4289.7998
$ws.Range("N132").Value = -10416.9998
